$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.693.60'
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").Value = '3.689.42'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'670.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").Value = "'160.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.19%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.43%  '
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").Value = "'0.442"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.75%  '
$ws.Range("D12").Value = "'0.0000233"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.27%  '
$ws.Range("D13").Value = "'33.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.44%  '
$ws.Range("D14").Value = '3.666.18'
$ws.Range("E14").Value = '  -0.31%  '
$ws.Range("D15").Value = '69.653.67'
$ws.Range("E15").Value = '  +0.67%  '
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = "'16.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.82%  '
$ws.Range("D18").Value = "'6.50"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.55%  '
$ws.Range("D19").Value = "'471.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = "'9.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.41%  '
$ws.Range("E21").Value = '  +0.28%  '
$ws.Range("D22").Value = "'79.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").Value = '3.837.38'
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").Value = "'0.0000127"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.78%  '
$ws.Range("E26").Value = '  +0.96%  '
$ws.Range("D27").Value = "'9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").Value = "'2.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").Value = "'2.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.54%  '
$ws.Range("E31").Value = '  +5.12%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = "'26.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = "'6.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.24%  '
$ws.Range("D35").Value = '3.686.96'
$ws.Range("E35").Value = '  +0.67%  '
$ws.Range("D36").Value = "'8.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.26%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = "'2.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.61%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").Value = "'176.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = "'0.0908"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.49%  '
$ws.Range("D43").Value = "'0.936"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("D44").Value = "'47.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("E45").Value = '  +3.20%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = "'1.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.75%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = "'27.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.65%  '
$ws.Range("E48").Value = '  -0.65%  '
$ws.Range("E49").Value = '  +0.93%  '
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("D51").Value = "'0.264"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.22%  '
